# Apply updated crypto price/volume data scraped on Wed Mar 22 22:56:09 UTC 2023.
# Column D ("Price") values are plain text in this sheet (not real numbers -
# e.g. "27.226.79" isn't a valid number at all, and precision such as the
# trailing zero in "13.50" must survive). A leading "'" forces Excel/COM to
# keep the assigned text as-is (quote-prefixed text) instead of silently
# re-parsing it into a Double and losing that formatting.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 44/45 coin identity + link also changed (list order swap): FraxShare <-> WEMIXTOKEN
$updates = @(
    @{ Row = 2; D = '27.254.16'; E = '  -3.13%  ' },
    @{ Row = 3; D = '1.731.53'; E = '  -3.70%  ' },
    @{ Row = 4; D = '''1.004'; E = '  -0.02%  ' },
    @{ Row = 5; D = '''322.66'; E = '  -3.97%  ' },
    @{ Row = 6; D = '''1.001'; E = '  +0.02%  ' },
    @{ Row = 7; D = '''0.4234'; E = '  -8.18%  ' },
    @{ Row = 8; D = '''0.3586'; E = '  -3.08%  ' },
    @{ Row = 9; D = '''44.86'; E = '  -0.65%  ' },
    @{ Row = 10; D = '''0.07389'; E = '  -3.48%  ' },
    @{ Row = 11; D = '''1.107'; E = '  -3.52%  ' },
    @{ Row = 12; E = '  +0.07%  ' },
    @{ Row = 13; D = '''21.51'; E = '  -4.44%  ' },
    @{ Row = 14; D = '''6.057'; E = '  -4.71%  ' },
    @{ Row = 15; D = '''7.138'; E = '  -3.45%  ' },
    @{ Row = 16; D = '1.730.47'; E = '  -3.57%  ' },
    @{ Row = 17; D = '''0.00001058' },
    @{ Row = 18; D = '''86.88' },
    @{ Row = 19; D = '''0.05971'; E = '  -11.16%  ' },
    @{ Row = 20; D = '''1.002'; E = '  +0.00%  ' },
    @{ Row = 21; D = '''16.83'; E = '  -3.35%  ' },
    @{ Row = 22; D = '''6.076'; E = '  -5.16%  ' },
    @{ Row = 23; D = '''0.5250'; E = '  -4.48%  ' },
    @{ Row = 24; D = '27.287.85'; E = '  -3.03%  ' },
    @{ Row = 25; D = '''11.31'; E = '  -4.95%  ' },
    @{ Row = 26; D = '''2.381'; E = '  -1.31%  ' },
    @{ Row = 27; D = '''20.06'; E = '  -3.25%  ' },
    @{ Row = 28; D = '''2.337'; E = '  -1.79%  ' },
    @{ Row = 29; D = '''148.27'; E = '  -2.44%  ' },
    @{ Row = 30; D = '1.925.97'; E = '  -3.77%  ' },
    @{ Row = 31; D = '''126.10'; E = '  -5.91%  ' },
    @{ Row = 32; D = '''1.185'; E = '  -6.07%  ' },
    @{ Row = 33; D = '''0.09077'; E = '  -5.46%  ' },
    @{ Row = 34; D = '''5.589'; E = '  -4.98%  ' },
    @{ Row = 35; D = '''3.557'; E = '  -12.19%  ' },
    @{ Row = 36; D = '''12.71'; E = '  +4.48%  ' },
    @{ Row = 37; D = '''0.2139'; E = '  -3.86%  ' },
    @{ Row = 38; D = '''5.065'; E = '  -4.19%  ' },
    @{ Row = 39; D = '''0.06050'; E = '  -4.80%  ' },
    @{ Row = 40; D = '''0.02239' },
    @{ Row = 41; D = '''0.6330'; E = '  -5.86%  ' },
    @{ Row = 42; D = '''1.186'; E = '  -4.12%  ' },
    @{ Row = 43; D = '''1.001'; E = '  +0.07%  ' },
    @{ Row = 44; B = 'FraxShare'; C = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'; D = '''7.892'; E = '  -2.50%  ' },
    @{ Row = 45; B = 'WEMIXTOKEN'; C = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'; D = '''1.404'; E = '  -7.43%  ' },
    @{ Row = 46; D = '''13.50'; E = '  -4.07%  ' },
    @{ Row = 47; D = '''3.714'; E = '  -3.34%  ' },
    @{ Row = 48; D = '''0.5806'; E = '  -5.81%  ' },
    @{ Row = 49; D = '''124.13'; E = '  -4.81%  ' },
    @{ Row = 50; D = '''1.943'; E = '  -5.61%  ' },
    @{ Row = 51; D = '''0.06808'; E = '  -4.52%  ' }
)

foreach ($u in $updates) {
    if ($u.ContainsKey("B")) { $ws.Range("B" + $u.Row).Value = $u.B }
    if ($u.ContainsKey("C")) { $ws.Range("C" + $u.Row).Value = $u.C }
    if ($u.ContainsKey("D")) { $ws.Range("D" + $u.Row).Value = $u.D }
    if ($u.ContainsKey("E")) { $ws.Range("E" + $u.Row).Value = $u.E }
}
